# Informal_road_avg_dropped_input.xlsx — "Add files via upload" re-save.
#
# The underlying data (headers + 85 rows of GRIP road-type / material
# intensity figures, including the AVERAGE()/cell-reference formulas in
# rows 7-13, 19-25, 31-37, 43-49, 55-61, 67-73, 79-85) is untouched by this
# commit - only sheet/view metadata changed:
#   * the sheet was renamed from the default "Sheet1" to
#     "Material Intensities Rousseau"
#   * the saved view no longer shows a scrolled-down window
#     (topLeftCell="A58") with B86 selected - it now opens scrolled back to
#     the top with A1 the active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the (only) worksheet.
$ws.Name = "Material Intensities Rousseau"

# Scroll the window back to the top-left corner and select A1, clearing the
# stale "topLeftCell=A58 / selection=A86:B86" view state that was left over
# from the previous save.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
